# Update the stack trace text in the red "failed:" run to reflect the
# move from AQL 7.x to 8.x (new stack frames / line numbers, and the
# redundant "/ by zero" summary line removed from the top of the trace).
$d = $word.ActiveDocument

$old = "    <---divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:`n`t/ by zero`njava.lang.ArithmeticException: / by zero`n`tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:77)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:569)`n`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:186)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:659)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2247)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:340)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:324)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:912)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:497)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)`n`tat jdk.internal.reflect.GeneratedMethodAccessor7.invoke(Unknown Source)`n`tat java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.base/java.lang.reflect.Method.invoke(Method.java:569)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)`n"
$new = "    <---divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:`n`tjava.lang.ArithmeticException: / by zero`n`t`tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)`n`t`tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)`n`t`tat java.base/java.lang.reflect.Method.invoke(Method.java:580)`n`t`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:170)`n`t`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:231)`n`t`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallService(EvaluationServices.java:122)`n`t`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCall(EvaluationServices.java:237)`n`t`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:170)`n`t`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:244)`n`t`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`t`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:146)`n`t`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:53)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:674)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)`n`t`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`t`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:350)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`t`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`t`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)`n`t`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:334)`n`t`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:878)`n`t`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:499)`n`t`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)`n`t`tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)`n`t`tat java.base/java.lang.reflect.Method.invoke(Method.java:580)`n`t`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)`n`t`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`t`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)`n`t`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`t`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`t`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`t`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`t`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`t`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`t`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`t`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`t`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`t`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`t`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`t`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`t`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`t`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`t`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`t`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`t`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`t`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`t`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`t`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`t`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`t`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`t`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`t`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`t`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)`n`t`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)`n`t`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)`n`t`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)`n`t`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)`n`t`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)`n`t"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    Write-Output "WARNING: stack trace text was not found/replaced"
} else {
    Write-Output "Stack trace text replaced successfully"
}
